# Add a "Job Posting" header row + one data row (Job_Id = JD_001) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$headers = @(
    "Job_Id",
    "Job_Title",
    "Job_Description",
    "Total_Years_Min_Exp",
    "Total_Years_Max_Exp",
    "Work_Mode",
    "Job_Location",
    "LinkedIn_Poster",
    "LinkedIn_Posted",
    "Resume_received",
    "Resume_downloaded"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Build the header formatting (bold font, thin box border, centered
# horizontally, top-aligned vertically) once on a scratch cell, then copy
# that single combined format onto the whole header range in one shot so we
# don't leave a trail of partial intermediate cell styles behind.
$tmpl = $ws.Cells.Item(100, 100)
$tmpl.Font.Bold = $true
$tmpl.HorizontalAlignment = -4108   # xlHAlignCenter
$tmpl.VerticalAlignment = -4160     # xlVAlignTop
$tmpl.Borders.LineStyle = 1         # xlContinuous (thin)

$tmpl.Copy()
$ws.Range("A1:K1").PasteSpecial(-4122)   # xlPasteFormats
$tmpl.Clear()
$excel.CutCopyMode = $false

# ---- Data row (row 2) ----
$description = "We are seeking a Junior RPA Developer to design, develop, and support automation solutions.`nCollaborate with teams to streamline business processes using RPA tools like UiPath or Automation Anywhere. Join Akkodis to grow your skills in a dynamic, tech-driven environment"

$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "Junior RPA Developer"
$ws.Range("C2").Value = $description
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = "Remote"
$ws.Range("G2").Value = "Hyderabad, Telangana, India"

# Reset row 2's height back to the default (no explicit custom height),
# since assigning the multi-line description auto-expands it.
$ws.Rows.Item(2).AutoFit()
